$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and two pairs of row swaps: rows 12/13 and 48/49)

# Row 2
$ws.Range("D2").Value = "'28.257.51"
$ws.Range("E2").Value = '  +3.09%  '

# Row 3
$ws.Range("D3").Value = "'1.814.76"
$ws.Range("E3").Value = '  +4.09%  '

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").Value = "'329.88"
$ws.Range("E5").Value = '  +2.39%  '

# Row 6
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = '  -0.12%  '

# Row 7
$ws.Range("D7").Value = "'0.4393"
$ws.Range("E7").Value = '  +4.19%  '

# Row 8
$ws.Range("D8").Value = "'0.3705"
$ws.Range("E8").Value = '  +3.07%  '

# Row 9
$ws.Range("D9").Value = "'44.93"
$ws.Range("E9").Value = '  -0.33%  '

# Row 10
$ws.Range("D10").Value = "'0.07734"
$ws.Range("E10").Value = '  +4.37%  '

# Row 11
$ws.Range("D11").Value = "'1.134"
$ws.Range("E11").Value = '  +2.12%  '

# Row 12
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = '  -0.09%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = "'22.15"
$ws.Range("E13").Value = '  +3.17%  '

# Row 14
$ws.Range("D14").Value = "'6.315"
$ws.Range("E14").Value = '  +3.90%  '

# Row 15
$ws.Range("D15").Value = "'7.586"
$ws.Range("E15").Value = '  +5.77%  '

# Row 16
$ws.Range("D16").Value = "'1.832.04"
$ws.Range("E16").Value = '  +4.86%  '

# Row 17
$ws.Range("D17").Value = "'92.96"
$ws.Range("E17").Value = '  +7.38%  '

# Row 18
$ws.Range("D18").Value = "'0.00001085"
$ws.Range("E18").Value = '  +1.92%  '

# Row 19
$ws.Range("D19").Value = "'0.06519"
$ws.Range("E19").Value = '  +8.04%  '

# Row 20
$ws.Range("D20").Value = "'0.9993"
$ws.Range("E20").Value = '  -0.15%  '

# Row 21
$ws.Range("D21").Value = "'17.57"
$ws.Range("E21").Value = '  +4.18%  '

# Row 22
$ws.Range("D22").Value = "'6.242"
$ws.Range("E22").Value = '  +2.47%  '

# Row 23
$ws.Range("D23").Value = "'28.294.53"
$ws.Range("E23").Value = '  +3.10%  '

# Row 24
$ws.Range("D24").Value = "'11.68"
$ws.Range("E24").Value = '  +2.70%  '

# Row 25
$ws.Range("D25").Value = "'2.149"
$ws.Range("E25").Value = '  -9.54%  '

# Row 26
$ws.Range("D26").Value = "'20.87"
$ws.Range("E26").Value = '  +2.14%  '

# Row 27
$ws.Range("D27").Value = "'157.62"
$ws.Range("E27").Value = '  +5.11%  '

# Row 28
$ws.Range("D28").Value = "'2.034.50"
$ws.Range("E28").Value = '  +4.78%  '

# Row 29
$ws.Range("D29").Value = "'2.298"
$ws.Range("E29").Value = '  -4.14%  '

# Row 30
$ws.Range("D30").Value = "'129.06"
$ws.Range("E30").Value = '  +2.02%  '

# Row 31
$ws.Range("D31").Value = "'1.206"
$ws.Range("E31").Value = '  +2.32%  '

# Row 32
$ws.Range("D32").Value = "'5.914"
$ws.Range("E32").Value = '  +3.95%  '

# Row 33
$ws.Range("D33").Value = "'0.09223"
$ws.Range("E33").Value = '  +1.53%  '

# Row 34
$ws.Range("D34").Value = "'3.621"
$ws.Range("E34").Value = '  +0.37%  '

# Row 35
$ws.Range("D35").Value = "'13.06"
$ws.Range("E35").Value = '  +1.35%  '

# Row 36
$ws.Range("D36").Value = "'0.02362"
$ws.Range("E36").Value = '  +4.02%  '

# Row 37
$ws.Range("D37").Value = "'0.2186"
$ws.Range("E37").Value = '  +2.65%  '

# Row 38
$ws.Range("D38").Value = "'5.186"
$ws.Range("E38").Value = '  +2.57%  '

# Row 39
$ws.Range("D39").Value = "'0.6589"
$ws.Range("E39").Value = '  +3.42%  '

# Row 40
$ws.Range("D40").Value = "'0.06221"
$ws.Range("E40").Value = '  +2.95%  '

# Row 41
$ws.Range("D41").Value = "'1.198"
$ws.Range("E41").Value = '  +0.83%  '

# Row 42
$ws.Range("D42").Value = "'8.153"
$ws.Range("E42").Value = '  +2.26%  '

# Row 43
$ws.Range("D43").Value = "'1.422"
$ws.Range("E43").Value = '  -0.40%  '

# Row 44
$ws.Range("D44").Value = "'0.9992"

# Row 45
$ws.Range("D45").Value = "'13.96"
$ws.Range("E45").Value = '  +1.65%  '

# Row 46
$ws.Range("D46").Value = "'0.6118"
$ws.Range("E46").Value = '  +5.16%  '

# Row 47
$ws.Range("D47").Value = "'3.770"
$ws.Range("E47").Value = '  +1.57%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = "'127.04"
$ws.Range("E48").Value = '  +1.75%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'2.038"
$ws.Range("E49").Value = '  +4.36%  '

# Row 50
$ws.Range("D50").Value = "'1.159"
$ws.Range("E50").Value = '  +5.38%  '

# Row 51
$ws.Range("D51").Value = "'0.07023"
$ws.Range("E51").Value = '  +2.61%  '
